# Hortaliza, Vega Modelo de Temuco - Achicoria: add two new weekly price
# rows into the historical log. The sheet keeps one row per market
# observation ordered roughly by date; two new observations are inserted
# (one becomes row 14, one becomes row 18), pushing the existing rows
# below them down and extending the used range from A1:R37 to A1:R39.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert the two new rows at their target positions. Inserting at row 14
# first (shifting old rows 14-37 down to 15-38), then inserting again at
# row 18 (shifting the now-current rows 18-38 down to 19-39) reproduces
# the final layout where the untouched historical rows keep their
# relative order around the two newcomers.
$ws.Rows("14").Insert()
$ws.Rows("18").Insert()

# New row 14: 2022-06-22 (serial 44720), Vega Modelo de Temuco observation.
$ws.Cells.Item(14, 1).Value = 10
$ws.Cells.Item(14, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(14, 3).Value = "La Araucanía"
$ws.Cells.Item(14, 4).Value = 44720
$ws.Cells.Item(14, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(14, 5).Value = 9
$ws.Cells.Item(14, 6).Value = 100112010
$ws.Cells.Item(14, 7).Value = "Achicoria"
$ws.Cells.Item(14, 8).Value = "Sin especificar"
$ws.Cells.Item(14, 9).Value = "Primera"
$ws.Cells.Item(14, 10).Value = 100
$ws.Cells.Item(14, 11).Value = 10000
$ws.Cells.Item(14, 12).Value = 10000
$ws.Cells.Item(14, 13).Value = 10000
$ws.Cells.Item(14, 14).Value = "$/caja 18 unidades"
$ws.Cells.Item(14, 15).Value = "Región Metropolitana"
$ws.Cells.Item(14, 16).Value = 556
$ws.Cells.Item(14, 17).Value = 18
$ws.Cells.Item(14, 18).Value = "Hortaliza"

# New row 18: 2022-06-23 (serial 44721), Vega Modelo de Temuco observation.
$ws.Cells.Item(18, 1).Value = 10
$ws.Cells.Item(18, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(18, 3).Value = "La Araucanía"
$ws.Cells.Item(18, 4).Value = 44721
$ws.Cells.Item(18, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(18, 5).Value = 9
$ws.Cells.Item(18, 6).Value = 100112010
$ws.Cells.Item(18, 7).Value = "Achicoria"
$ws.Cells.Item(18, 8).Value = "Sin especificar"
$ws.Cells.Item(18, 9).Value = "Primera"
$ws.Cells.Item(18, 10).Value = 80
$ws.Cells.Item(18, 11).Value = 10000
$ws.Cells.Item(18, 12).Value = 10000
$ws.Cells.Item(18, 13).Value = 10000
$ws.Cells.Item(18, 14).Value = "$/caja 18 unidades"
$ws.Cells.Item(18, 15).Value = "Región Metropolitana"
$ws.Cells.Item(18, 16).Value = 556
$ws.Cells.Item(18, 17).Value = 18
$ws.Cells.Item(18, 18).Value = "Hortaliza"
